$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Exportar Planilha" sheet: append row 45 with the new monthly figures
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Exportar Planilha")

# Columns A/B hold text-looking-like-numbers ("2017", "201708"), exactly like
# the existing rows above them. A plain .Value assignment would be silently
# coerced to a numeric value by Excel (since the cells still carry the
# default "General" format), so instead enter short text formulas and
# immediately freeze them down to static values -- that keeps the result as
# genuine text without leaving any NumberFormat override on the cell.
$ws.Range("A45").Formula = '="2017"'
$ws.Range("B45").Formula = '="201708"'
$ws.Range("A45:B45").Copy()
$ws.Range("A45:B45").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Columns C:R reuse the same cell style as the row above (right-aligned,
# "Dialog" font) -- copy that formatting down before writing the values.
$ws.Range("C44:R44").Copy()
$ws.Range("C45:R45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C45").Value = 10197385.38
$ws.Range("D45").Value = 79.02
$ws.Range("E45").Value = 271.15
$ws.Range("F45").Value = 139.21
$ws.Range("G45").Value = 48.3
$ws.Range("H45").Value = 514417.0
$ws.Range("I45").Value = 3.99
$ws.Range("J45").Value = 13.68
$ws.Range("K45").Value = 7.02
$ws.Range("L45").Value = 2.44
$ws.Range("M45").Value = 19.823188930381384
$ws.Range("N45").Value = 129040.0
$ws.Range("O45").Value = 37608.0
$ws.Range("P45").Value = 73252.0
$ws.Range("Q45").Value = 211147.0
$ws.Range("R45").Value = 0.29144451332920024

# ---------------------------------------------------------------------------
# 2) "SQL" sheet: widen the date window to 01/08/2017, comment out / add a
#    few alternate filter lines for age bracket and beneficiary age range.
# ---------------------------------------------------------------------------
$sql = $wb.Worksheets.Item("SQL")

$newSql = "select  substr(fc.ID_TEMPO_MES_ANO_REF,1,4) ano,`n        substr(fc.ID_TEMPO_MES_ANO_REF,1,6) anomes,`n        sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0))  VT,        `n        round((sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)))/(tbt.BT),2) VM,`n        round((sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)))/(count(distinct fc.COD_TS)),2) VA,`n        round((sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)))/(count(distinct fc.CONTA)),2) VC,`n        round((sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)))/(count(1)),2) VI,`n        sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0))  QP,        `n        round((sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)))/(tbt.BT),2) QM,        `n        round((sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)))/(count(distinct fc.COD_TS)),2) QA,`n        round((sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)))/(count(distinct fc.CONTA)),2) QC,`n        round((sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)))/(count(1)),2) QI,`n        (sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)))/(sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0))) VQ,`n        tbt.BT,`n        count(distinct fc.COD_TS) BA,`n        count(distinct fc.CONTA) CT,`n        count(1) IT,`n        count(distinct fc.COD_TS)/tbt.BT PA        `nfrom    TS.FAT_ITEM_CONTA fc,`n        (select  to_char(pcm.mes_ano_ref,'RRRRMM') anomes,`n                sum(qtd_ativos) BT`n        from    ts.posicao_cadastro_mes pcm           `n        where   pcm.mes_ano_ref between to_date ('01/01/2014','dd/mm/yyyy') and to_date ('01/08/2017','dd/mm/yyyy')`n        --and     pcm.COD_FAIXA_ETARIA = 10`n        and     pcm.COD_FAIXA_ETARIA = 1`n        --and     pcm.COD_FAIXA_ETARIA not in (1, 10)        `n        group by   `n                to_char(pcm.mes_ano_ref,'RRRRMM')`n        order by 1) tbt,`n        /**/(select  b.COD_TS,`n                b.NOME_ASSOCIADO,`n                b.DATA_NASCIMENTO,`n                2017-to_number(to_char(b.DATA_NASCIMENTO,'RRRR')) idade`n        from    ts.dim_beneficiario b`n        --where   2017-to_number(to_char(b.DATA_NASCIMENTO,'RRRR')) > 58) bnf`n        where   2017-to_number(to_char(b.DATA_NASCIMENTO,'RRRR')) < 19) bnf`n        --where   2017-to_number(to_char(b.DATA_NASCIMENTO,'RRRR')) between 19 and 58) bnf        `nwhere   substr(fc.ID_TEMPO_MES_ANO_REF,1,6) = tbt.anomes`nand     fc.COD_TS = bnf.COD_TS`n--and     substr(fc.ID_TEMPO_MES_ANO_REF,1,6) in ('201605')--, '201606', '201705', '201706')`ngroup by substr(fc.ID_TEMPO_MES_ANO_REF,1,4),`n         substr(fc.ID_TEMPO_MES_ANO_REF,1,6),`n         tbt.BT`nhaving  sum(nvl(fc.QTD_ITEM,0)) - sum(nvl(fc.QTD_GLOSADO,0)) > 0`nor      sum(nvl(fc.VAL_APROVADO_ITEM,0)) + sum(nvl(fc.VALOR_PAGO_REVISAO,0)) > 0`norder by 2"

$sql.Range("A2").Value = $newSql
